$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.215935230255127
$ws.Range("B1").Value = 2.293593645095825
$ws.Range("C1").Value = 2.956059217453003
$ws.Range("D1").Value = 3.586608648300171
$ws.Range("E1").Value = 1.585543155670166
